# Add a new daily case-count column (R) to the "CodeBook" sheet, mirroring
# the existing D:Q daily columns, and extend the row-25 SUM() totals to
# cover it. Finally, move the sheet's active selection to the newly
# completed cells (matches the author re-saving after typing the new
# column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeBook")

# New case counts for 2020-04-14 (column R), one per county/city row (3-24).
$newCounts = @{
    3  = 112
    4  = 37
    5  = 14
    6  = 29
    7  = 6
    8  = 11
    9  = 3
    10 = 88
    11 = 44
    12 = 6
    13 = 2
    14 = 3
    15 = 18
    16 = 2
    17 = 4
    18 = 2
    19 = 7
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
}

foreach ($row in $newCounts.Keys) {
    $ws.Cells.Item($row, 18).Value = $newCounts[$row]
}

# Copy the formatting (style) of column Q down into the new column R so the
# new cells look like the rest of the data block.
$ws.Range("Q3:Q24").Copy()
$ws.Range("R3:R24").PasteSpecial(-4122)

# Extend the running total row to include the new column.
$ws.Range("R25").Formula = "=SUM(R3:R24)"
$ws.Range("Q25").Copy()
$ws.Range("R25").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the selection to the newly-added total cells, as in the saved file.
$ws.Activate()
$ws.Range("Q25:R25").Select()
